# Edit: 2020-03-22 commit
# 1) Table on slide 6 gets a new table style id.
# 2) The deck's theme (theme1.xml, used by the slide master/"Integral" design)
#    is repainted with the standard Office Theme colour palette (the palette
#    that used to live in theme2.xml, the notes-master theme).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{5B4322AA-88EA-46A5-AFD3-8AC6DD8388FF}", $true)
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Map of ThemeColorSchemeIndex (1..12) -> target RGB (0xRRGGBB) for the
# "Office Theme" colour scheme.
$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgbHex = $officeColors[$i - 1]
    $r = ($rgbHex -shr 16) -band 0xFF
    $g = ($rgbHex -shr 8) -band 0xFF
    $b = $rgbHex -band 0xFF
    $bgr = ($b -shl 16) -bor ($g -shl 8) -bor $r
    $tcs.Colors($i).RGB = $bgr
}
